$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the two requirement lines so the LOT2028 (weak requisite) line now comes
# before the LOT2038 (set indication) line, matching the reordering of the
# corresponding <si> entries in sharedStrings.xml.
$lot2038 = $ws.Range("B23").Value()
$lot2028 = $ws.Range("B24").Value()

$ws.Range("B23").Value = $lot2028
$ws.Range("C23").Value = $lot2028
$ws.Range("B24").Value = $lot2038
$ws.Range("C24").Value = $lot2038
